$wb = $excel.ActiveWorkbook

$wsDBD = $wb.Worksheets.Item("DBD")
$wsDBS = $wb.Worksheets.Item("DBS")

# Rename field "BonusNo" -> "LogNo" (Key ID row, field row, Index1 row)
$wsDBD.Range("C3").Value = "LogNo"
$wsDBD.Range("C5").Value = "WorkMonth,LogNo"
$wsDBD.Range("B9").Value = "LogNo"

# Field type change DATE -> TIMESTAMP for CreateDate / LastUpdate rows
$wsDBD.Range("D28").Value = "TIMESTAMP"
$wsDBD.Range("D30").Value = "TIMESTAMP"

# Update view state: active sheet DBD, scrolled down, selection at D28
$wsDBD.Activate()
$excel.ActiveWindow.ScrollRow = 25
$wsDBD.Range("D28").Select()

$wsDBS.Range("B2").Select()
